# Form the consolidated report: populate the "Absent" (column H) values
# that were left blank/uncalculated for rows 4, 7, 9, 10, 11 and 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Absent becomes 1 (was 0)
$ws.Range("H4").Value = 1

# Row 7: Absent gets filled in as numeric 0 (was an empty/inline string cell)
$ws.Range("H7").Value = 0

# Row 9: Absent becomes 1 (was 0)
$ws.Range("H9").Value = 1

# Row 10: Absent gets filled in as numeric 0 (was an empty/inline string cell)
$ws.Range("H10").Value = 0

# Row 11: Absent becomes 1 (was 0)
$ws.Range("H11").Value = 1

# Row 12: Absent gets filled in as numeric 0 (was an empty/inline string cell)
$ws.Range("H12").Value = 0
